$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 185
$ws.Range("A5").Value = 44
$ws.Range("A6").Value = 45
$ws.Range("A8").Value = 100
$ws.Range("A11").Value = 96
$ws.Range("A14").Value = 190
$ws.Range("A18").Value = 171
$ws.Range("A19").Value = 175
$ws.Range("A20").Value = 4
$ws.Range("A21").Value = 3
$ws.Range("A22").Value = 49
$ws.Range("A24").Value = 67
